# Rotate the "Fecha"/Volumen/Precio/Unidad/Origen/Precio-Kg/Kg-unidad data
# across rows 3, 4 and 5 (row 5 -> row 3, row 3 -> row 4, row 4 -> row 5),
# keeping the market/region/product descriptive columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that move.
$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

$row3 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Range("$col" + "3").Value2
    $row4[$col] = $ws.Range("$col" + "4").Value2
    $row5[$col] = $ws.Range("$col" + "5").Value2
}

# Apply the rotation: new row3 = old row5, new row4 = old row3, new row5 = old row4.
foreach ($col in $cols) {
    $ws.Range("$col" + "3").Value = $row5[$col]
    $ws.Range("$col" + "4").Value = $row3[$col]
    $ws.Range("$col" + "5").Value = $row4[$col]
}
